$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Anonymize invoice content (only cells whose text actually changes) ---
$ws.Range("A7").Value = "Paslaugos teikėjas: Vardas Pavardė"
$ws.Range("A8").Value = "Individualios veiklos vykdymo pažyma Nr. 1234567"
$ws.Range("A9").Value = "Adresas: Gatvės g. 1, Miestas"
$ws.Range("A10").Value = "Banko sąskaita: LT123456789123456789"
$ws.Range("A11").Value = "Bankas: ManoBankas"
$ws.Range("A12").Value = "Tel. Nr.: +370 61234567"
$ws.Range("A13").Value = "El. Paštas:  paštas@gmail.com"
$ws.Range("B16").Value = "Paslauga"
$ws.Range("A20").Value = "Pastabos: "
$ws.Range("A25").Value = "Vadovė: Vardas Pavardė"

# --- View / selection changes ---
$excel.ActiveWindow.Zoom = 115
$ws.Range("B24").Select()
